$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: "Average Sales"/"Minimum Sales"/"Max Sales"/"Total Sales" shift and two
# new headers are introduced, pushing the old ones over.
$ws.Range("B1").Value = "Monthly sales"
$ws.Range("C1").Value = "% of change sales"
$ws.Range("D1").Value = "Minimum Sales"
$ws.Range("E1").Value = "Max Sales"

# "Total Sales" and "Average Sales" move out of the header row into new summary rows.
$ws.Range("A17").Value = "Total Sales"
$ws.Range("A18").Value = "Average Sales"
